$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 851
$ws1.Range("F23").Value = 1109
$ws1.Range("F24").Value = 2028
$ws1.Range("F28").Value = 51
$ws1.Range("F29").Value = 2120
$ws1.Range("F30").Value = 78

# Sheet "全部类型" (all types) - same updates, rows shifted by +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 851
$ws4.Range("F24").Value = 1109
$ws4.Range("F25").Value = 2028
$ws4.Range("F29").Value = 51
$ws4.Range("F30").Value = 2120
$ws4.Range("F31").Value = 78
